$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 7694807.5  # H43: 3849028.2 -> 7694807.5
$ws.Cells.Item(43, 9).Value = 7694807.5  # I43: 5130871.5 -> 7694807.5
$ws.Cells.Item(43, 10).Value = 0  # J43: 3498 -> 0
$ws.Cells.Item(43, 11).Value = 7694807.5  # K43: 5130871.5 -> 7694807.5
$ws.Cells.Item(43, 12).Value = 0  # L43: 3498 -> 0
$ws.Cells.Item(43, 13).Value = -7694738.5  # M43: -5130802.5 -> -7694738.5
$ws.Cells.Item(43, 14).ClearContents()  # N43: -3636 -> (deleted)
$ws.Cells.Item(86, 8).Value = 1156969  # H86: 1117107.5 -> 1156969
$ws.Cells.Item(86, 9).Value = 1797773  # I86: 1703205.1 -> 1797773
$ws.Cells.Item(86, 11).Value = 1797773  # K86: 1703205.1 -> 1797773
$ws.Cells.Item(86, 13).Value = -1796650  # M86: -1702082.1 -> -1796650
$ws.Cells.Item(89, 8).Value = 1156969  # H89: 1117107.5 -> 1156969
$ws.Cells.Item(89, 9).Value = 1797773  # I89: 1703205.1 -> 1797773
$ws.Cells.Item(89, 11).Value = 8988865  # K89: 8516025.5 -> 8988865
$ws.Cells.Item(89, 13).Value = -8983249  # M89: -8510409.5 -> -8983249
$ws.Cells.Item(92, 8).Value = 635.9167  # H92: 603.9231 -> 635.9167
$ws.Cells.Item(92, 9).Value = 520.5454999999999  # I92: 495.5 -> 520.5454999999999
$ws.Cells.Item(92, 11).Value = 520.5454999999999  # K92: 495.5 -> 520.5454999999999
$ws.Cells.Item(92, 13).Value = 727.4545000000001  # M92: 752.5 -> 727.4545000000001
$ws.Cells.Item(121, 8).Value = 3877.4827  # H121: 3840 -> 3877.4827
$ws.Cells.Item(121, 10).Value = 3877.4827  # J121: 3840 -> 3877.4827
$ws.Cells.Item(121, 12).Value = 11632.4481  # L121: 11520 -> 11632.4481
$ws.Cells.Item(121, 14).Value = -15126.4481  # N121: -15014 -> -15126.4481
$ws.Cells.Item(129, 8).Value = 2649.4614  # H129: 2362.9412 -> 2649.4614
$ws.Cells.Item(129, 9).Value = 694  # I129: 731.9 -> 694
$ws.Cells.Item(129, 10).Value = 5778.2  # J129: 4693 -> 5778.2
$ws.Cells.Item(129, 11).Value = 2082  # K129: 2195.7 -> 2082
$ws.Cells.Item(129, 12).Value = 17334.6  # L129: 14079 -> 17334.6
$ws.Cells.Item(129, 13).Value = 2918  # M129: 2804.3 -> 2918
$ws.Cells.Item(129, 14).Value = -27334.6  # N129: -24079 -> -27334.6
$ws.Cells.Item(138, 8).Value = 3406.7068  # H138: 4211.054 -> 3406.7068
$ws.Cells.Item(138, 9).Value = 1402.375  # I138: 1353.6471 -> 1402.375
$ws.Cells.Item(138, 10).Value = 4170.2617  # J138: 6639.85 -> 4170.2617
$ws.Cells.Item(138, 11).Value = 4207.125  # K138: 4060.9413 -> 4207.125
$ws.Cells.Item(138, 12).Value = 12510.7851  # L138: 19919.55 -> 12510.7851
$ws.Cells.Item(138, 13).Value = 932.875  # M138: 1079.0587 -> 932.875
$ws.Cells.Item(138, 14).Value = -22790.7851  # N138: -30199.55 -> -22790.7851
$ws.Cells.Item(140, 8).Value = 57073.9  # H140: 58783.8 -> 57073.9
$ws.Cells.Item(140, 10).Value = 55558.89  # J140: 57458.777 -> 55558.89
$ws.Cells.Item(140, 12).Value = 55558.89  # L140: 57458.777 -> 55558.89
$ws.Cells.Item(140, 14).Value = -65918.89  # N140: -67818.777 -> -65918.89
$ws.Cells.Item(141, 8).Value = 5995  # H141: 6140.923 -> 5995
$ws.Cells.Item(141, 9).Value = 5274.88  # I141: 5402.9585 -> 5274.88
$ws.Cells.Item(141, 11).Value = 15824.64  # K141: 16208.8755 -> 15824.64
$ws.Cells.Item(141, 13).Value = -10644.64  # M141: -11028.8755 -> -10644.64

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 3000  # H3: 4000 -> 3000
$ws.Cells.Item(3, 10).Value = 0  # J3: 5000 -> 0
$ws.Cells.Item(3, 12).Value = 0  # L3: 5000 -> 0
$ws.Cells.Item(3, 14).ClearContents()  # N3: -5230 -> (deleted)
$ws.Cells.Item(4, 8).Value = 242.5  # H4: 371.33334 -> 242.5
$ws.Cells.Item(4, 9).Value = 242.5  # I4: 450 -> 242.5
$ws.Cells.Item(4, 10).Value = 0  # J4: 214 -> 0
$ws.Cells.Item(4, 11).Value = 242.5  # K4: 450 -> 242.5
$ws.Cells.Item(4, 12).Value = 0  # L4: 214 -> 0
$ws.Cells.Item(4, 13).Value = -126.5  # M4: -334 -> -126.5
$ws.Cells.Item(4, 14).ClearContents()  # N4: -446 -> (deleted)
$ws.Cells.Item(16, 8).Value = 404.16666  # H16: 0 -> 404.16666
$ws.Cells.Item(16, 9).Value = 404.16666  # I16: 0 -> 404.16666
$ws.Cells.Item(16, 11).Value = 404.16666  # K16: 0 -> 404.16666
$ws.Cells.Item(16, 13).Value = -117.16666  # M16: None -> -117.16666
$ws.Cells.Item(32, 8).Value = 4364.3335  # H32: 4438.1343 -> 4364.3335
$ws.Cells.Item(32, 9).Value = 2347  # I32: 2395.2642 -> 2347
$ws.Cells.Item(32, 10).Value = 11626.733  # J32: 12171.857 -> 11626.733
$ws.Cells.Item(32, 11).Value = 2347  # K32: 2395.2642 -> 2347
$ws.Cells.Item(32, 12).Value = 11626.733  # L32: 12171.857 -> 11626.733
$ws.Cells.Item(32, 13).Value = -2060  # M32: -2108.2642 -> -2060
$ws.Cells.Item(32, 14).Value = -12200.733  # N32: -12745.857 -> -12200.733
$ws.Cells.Item(140, 8).Value = 111809.664  # H140: 101156 -> 111809.664
$ws.Cells.Item(140, 9).Value = 0  # I140: 80390 -> 0
$ws.Cells.Item(140, 10).Value = 111809.664  # J140: 115000 -> 111809.664
$ws.Cells.Item(140, 11).Value = 0  # K140: 80390 -> 0
$ws.Cells.Item(140, 12).Value = 111809.664  # L140: 115000 -> 111809.664
$ws.Cells.Item(140, 13).ClearContents()  # M140: -75210 -> (deleted)
$ws.Cells.Item(140, 14).Value = -122169.664  # N140: -125360 -> -122169.664
$ws.Cells.Item(141, 8).Value = 119999.5  # H141: 120000 -> 119999.5
$ws.Cells.Item(141, 9).Value = 109999  # I141: 100000 -> 109999
$ws.Cells.Item(141, 11).Value = 109999  # K141: 100000 -> 109999
$ws.Cells.Item(141, 13).Value = -104819  # M141: -94820 -> -104819

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 625.75  # H22: 619.6 -> 625.75
$ws.Cells.Item(22, 9).Value = 586.7143  # I22: 549.625 -> 586.7143
$ws.Cells.Item(22, 10).Value = 899  # J22: 899.5 -> 899
$ws.Cells.Item(22, 11).Value = 586.7143  # K22: 549.625 -> 586.7143
$ws.Cells.Item(22, 12).Value = 899  # L22: 899.5 -> 899
$ws.Cells.Item(22, 13).Value = -413.7143  # M22: -376.625 -> -413.7143
$ws.Cells.Item(22, 14).Value = -1245  # N22: -1245.5 -> -1245

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 2500  # H4: 0 -> 2500
$ws.Cells.Item(4, 9).Value = 2500  # I4: 0 -> 2500
$ws.Cells.Item(4, 11).Value = 2500  # K4: 0 -> 2500
$ws.Cells.Item(4, 13).Value = -2388  # M4: None -> -2388
$ws.Cells.Item(31, 8).Value = 3248.15  # H31: 3352.8333 -> 3248.15
$ws.Cells.Item(31, 9).Value = 1012.4286  # I31: 1027.8462 -> 1012.4286
$ws.Cells.Item(31, 10).Value = 8464.833000000001  # J31: 9397.799999999999 -> 8464.833000000001
$ws.Cells.Item(31, 11).Value = 1012.4286  # K31: 1027.8462 -> 1012.4286
$ws.Cells.Item(31, 12).Value = 8464.833000000001  # L31: 9397.799999999999 -> 8464.833000000001
$ws.Cells.Item(31, 13).Value = -717.4286  # M31: -732.8462 -> -717.4286
$ws.Cells.Item(31, 14).Value = -9054.833000000001  # N31: -9987.799999999999 -> -9054.833000000001
$ws.Cells.Item(34, 8).Value = 3248.15  # H34: 3352.8333 -> 3248.15
$ws.Cells.Item(34, 9).Value = 1012.4286  # I34: 1027.8462 -> 1012.4286
$ws.Cells.Item(34, 10).Value = 8464.833000000001  # J34: 9397.799999999999 -> 8464.833000000001
$ws.Cells.Item(34, 11).Value = 1012.4286  # K34: 1027.8462 -> 1012.4286
$ws.Cells.Item(34, 12).Value = 8464.833000000001  # L34: 9397.799999999999 -> 8464.833000000001
$ws.Cells.Item(34, 13).Value = -810.4286  # M34: -825.8462 -> -810.4286
$ws.Cells.Item(34, 14).Value = -8868.833000000001  # N34: -9801.799999999999 -> -8868.833000000001
$ws.Cells.Item(58, 8).Value = 717756.4  # H58: 772737.6 -> 717756.4
$ws.Cells.Item(58, 9).Value = 912254.4399999999  # I58: 1003179.9 -> 912254.4399999999
$ws.Cells.Item(58, 11).Value = 912254.4399999999  # K58: 1003179.9 -> 912254.4399999999
$ws.Cells.Item(58, 13).Value = -912051.4399999999  # M58: -1002976.9 -> -912051.4399999999
$ws.Cells.Item(62, 8).Value = 33150.1  # H62: 46343 -> 33150.1
$ws.Cells.Item(62, 9).Value = 3100  # I62: 5300 -> 3100
$ws.Cells.Item(62, 11).Value = 3100  # K62: 5300 -> 3100
$ws.Cells.Item(62, 13).Value = -2476  # M62: -4676 -> -2476
$ws.Cells.Item(65, 8).Value = 33150.1  # H65: 46343 -> 33150.1
$ws.Cells.Item(65, 9).Value = 3100  # I65: 5300 -> 3100
$ws.Cells.Item(65, 11).Value = 15500  # K65: 26500 -> 15500
$ws.Cells.Item(65, 13).Value = -12380  # M65: -23380 -> -12380
$ws.Cells.Item(132, 8).Value = 10765165  # H132: 7761117.5 -> 10765165
$ws.Cells.Item(132, 9).Value = 11914290  # I132: 8340201.5 -> 11914290
$ws.Cells.Item(132, 11).Value = 35742870  # K132: 25020604.5 -> 35742870
$ws.Cells.Item(132, 13).Value = -35740340  # M132: -25018074.5 -> -35740340
$ws.Cells.Item(136, 8).Value = 717756.4  # H136: 772737.6 -> 717756.4
$ws.Cells.Item(136, 9).Value = 912254.4399999999  # I136: 1003179.9 -> 912254.4399999999
$ws.Cells.Item(136, 11).Value = 2736763.32  # K136: 3009539.7 -> 2736763.32
$ws.Cells.Item(136, 13).Value = -2734213.32  # M136: -3006989.7 -> -2734213.32

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 3805469  # H4: 2114320.2 -> 3805469
$ws.Cells.Item(4, 9).Value = 1300311.2  # I4: 667010.9399999999 -> 1300311.2
$ws.Cells.Item(4, 10).Value = 13826100  # J4: 11521832 -> 13826100
$ws.Cells.Item(4, 11).Value = 3900933.6  # K4: 2001032.82 -> 3900933.6
$ws.Cells.Item(4, 12).Value = 41478300  # L4: 34565496 -> 41478300
$ws.Cells.Item(4, 13).Value = -3900821.6  # M4: -2000920.82 -> -3900821.6
$ws.Cells.Item(4, 14).Value = -41478524  # N4: -34565720 -> -41478524
$ws.Cells.Item(9, 8).Value = 526.94116  # H9: 482.05554 -> 526.94116
$ws.Cells.Item(9, 9).Value = 530.1667  # I9: 467.76923 -> 530.1667
$ws.Cells.Item(9, 11).Value = 1590.5001  # K9: 1403.30769 -> 1590.5001
$ws.Cells.Item(9, 13).Value = -1366.5001  # M9: -1179.30769 -> -1366.5001
$ws.Cells.Item(38, 8).Value = 204.41176  # H38: 194.77777 -> 204.41176
$ws.Cells.Item(38, 9).Value = 146.28572  # I38: 162.5 -> 146.28572
$ws.Cells.Item(38, 10).Value = 245.1  # J38: 210.91667 -> 245.1
$ws.Cells.Item(38, 11).Value = 438.85716  # K38: 487.5 -> 438.85716
$ws.Cells.Item(38, 12).Value = 735.3  # L38: 632.75001 -> 735.3
$ws.Cells.Item(38, 13).Value = -91.85716000000002  # M38: -140.5 -> -91.85716000000002
$ws.Cells.Item(38, 14).Value = -1429.3  # N38: -1326.75001 -> -1429.3
$ws.Cells.Item(63, 8).Value = 25000  # H63: 19000 -> 25000
$ws.Cells.Item(63, 9).Value = 0  # I63: 1000 -> 0
$ws.Cells.Item(63, 11).Value = 0  # K63: 3000 -> 0
$ws.Cells.Item(63, 13).ClearContents()  # M63: -2251 -> (deleted)
$ws.Cells.Item(66, 8).Value = 25000  # H66: 19000 -> 25000
$ws.Cells.Item(66, 9).Value = 0  # I66: 1000 -> 0
$ws.Cells.Item(66, 11).Value = 0  # K66: 9000 -> 0
$ws.Cells.Item(66, 13).ClearContents()  # M66: -5256 -> (deleted)
$ws.Cells.Item(101, 8).Value = 23888.666  # H101: 22599.6 -> 23888.666
$ws.Cells.Item(101, 10).Value = 23888.666  # J101: 22599.6 -> 23888.666
$ws.Cells.Item(101, 12).Value = 71665.99800000001  # L101: 67798.79999999999 -> 71665.99800000001
$ws.Cells.Item(101, 14).Value = -76533.99800000001  # N101: -72666.79999999999 -> -76533.99800000001
$ws.Cells.Item(113, 8).Value = 699.8  # H113: 671.2857 -> 699.8
$ws.Cells.Item(113, 9).Value = 798  # I113: 621.75 -> 798
$ws.Cells.Item(113, 10).Value = 675.25  # J113: 737.3333 -> 675.25
$ws.Cells.Item(113, 11).Value = 2394  # K113: 1865.25 -> 2394
$ws.Cells.Item(113, 12).Value = 2025.75  # L113: 2211.9999 -> 2025.75
$ws.Cells.Item(113, 13).Value = -224  # M113: 304.75 -> -224
$ws.Cells.Item(113, 14).Value = -6365.75  # N113: -6551.9999 -> -6365.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 130.07692  # H2: 125.18519 -> 130.07692
$ws.Cells.Item(2, 9).Value = 106  # I2: 99.86667 -> 106
$ws.Cells.Item(2, 10).Value = 158.16667  # J2: 156.83333 -> 158.16667
$ws.Cells.Item(2, 11).Value = 106  # K2: 99.86667 -> 106
$ws.Cells.Item(2, 12).Value = 158.16667  # L2: 156.83333 -> 158.16667
$ws.Cells.Item(2, 13).Value = 7  # M2: 13.13333 -> 7
$ws.Cells.Item(2, 14).Value = -384.16667  # N2: -382.83333 -> -384.16667
$ws.Cells.Item(46, 8).Value = 21996.5  # H46: 21996.666 -> 21996.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 501  # H22: 0 -> 501
$ws.Cells.Item(22, 9).Value = 501  # I22: 0 -> 501
$ws.Cells.Item(22, 11).Value = 501  # K22: 0 -> 501
$ws.Cells.Item(22, 13).Value = -206  # M22: None -> -206
$ws.Cells.Item(27, 8).Value = 501  # H27: 0 -> 501
$ws.Cells.Item(27, 9).Value = 501  # I27: 0 -> 501
$ws.Cells.Item(27, 11).Value = 501  # K27: 0 -> 501
$ws.Cells.Item(27, 13).Value = -394  # M27: None -> -394
$ws.Cells.Item(46, 8).Value = 7300.148  # H46: 7398.231 -> 7300.148
$ws.Cells.Item(46, 10).Value = 7285.16  # J46: 7390.7915 -> 7285.16
$ws.Cells.Item(46, 12).Value = 7285.16  # L46: 7390.7915 -> 7285.16
$ws.Cells.Item(46, 14).Value = -7661.16  # N46: -7766.7915 -> -7661.16
$ws.Cells.Item(93, 8).Value = 4655.6665  # H93: 4662 -> 4655.6665
$ws.Cells.Item(93, 9).Value = 980  # I93: 999 -> 980
$ws.Cells.Item(93, 11).Value = 980  # K93: 999 -> 980
$ws.Cells.Item(93, 13).Value = 268  # M93: 249 -> 268

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 60000  # H49: 65000 -> 60000
$ws.Cells.Item(49, 10).Value = 60000  # J49: 65000 -> 60000
$ws.Cells.Item(49, 12).Value = 60000  # L49: 65000 -> 60000
$ws.Cells.Item(49, 14).Value = -60460  # N49: -65460 -> -60460
$ws.Cells.Item(126, 8).Value = 1329.5555  # H126: 1413.8572 -> 1329.5555
$ws.Cells.Item(126, 9).Value = 1185.6666  # I126: 1220.5 -> 1185.6666
$ws.Cells.Item(126, 10).Value = 2049  # J126: 2574 -> 2049
$ws.Cells.Item(126, 11).Value = 3556.9998  # K126: 3661.5 -> 3556.9998
$ws.Cells.Item(126, 12).Value = 6147  # L126: 7722 -> 6147
$ws.Cells.Item(126, 13).Value = -1086.9998  # M126: -1191.5 -> -1086.9998
$ws.Cells.Item(126, 14).Value = -11087  # N126: -12662 -> -11087
$ws.Cells.Item(136, 8).Value = 9683.888999999999  # H136: 9817.223 -> 9683.888999999999
$ws.Cells.Item(136, 9).Value = 4541.8  # I136: 4749.263 -> 4541.8
$ws.Cells.Item(136, 10).Value = 12075.559  # J136: 12005.659 -> 12075.559
$ws.Cells.Item(136, 11).Value = 13625.4  # K136: 14247.789 -> 13625.4
$ws.Cells.Item(136, 12).Value = 36226.677  # L136: 36016.977 -> 36226.677
$ws.Cells.Item(136, 13).Value = -11075.4  # M136: -11697.789 -> -11075.4
$ws.Cells.Item(136, 14).Value = -41326.677  # N136: -41116.977 -> -41326.677
